$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concepts")

# Copy the "Display" (column C) value into the "Definition" (column D)
# column for every data row (rows 2-32).
for ($r = 2; $r -le 32; $r++) {
    $ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 3).Value2
}
